$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear manufacturer value in N2 (empty inline string cell)
$ws.Range("N2").Value = ""

# Update category text values in AL2 and AM2: replace " > " separator between
# GM Category and Subcategory with " | "
$ws.Range("AL2").Value = "Tovary a kategórie > GM Category | Subcategory"
$ws.Range("AM2").Value = "Tovary a kategórie > GM Category | Subcategory"
